# The deck currently carries the "Integral" theme on its (one and only)
# slide master (ppt/theme/theme1.xml) and the stock "Office Theme" on the
# notes master (ppt/theme/theme2.xml). The authored change swaps those two
# themes wholesale - the deck's slides end up themed with the plain
# "Office Theme" palette instead of "Integral".
#
# The font scheme (majorFont/minorFont) and the format scheme (fills,
# lines, effects) are already byte-for-byte identical between the two
# theme parts, so the only real difference is the 12-slot theme colour
# scheme (clrScheme). Re-pointing the deck at the "Office Theme" look is
# therefore just a matter of writing that theme's colour values into the
# live theme's ThemeColorScheme - PowerPoint has no VBA/COM surface for
# renaming the <a:theme>/<a:clrScheme> "name" attributes themselves, so
# those stay as-is (same as a real PowerPoint automation session would
# leave them).

$p = $ppt.ActivePresentation

function HexToRgbLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme, in ThemeColorScheme slot order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink.
# Slots 1-2 (dk1/lt1) are already 000000/FFFFFF in both themes, so only
# 3-12 actually need to change.
$officeColors = @{
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme

foreach ($slot in $officeColors.Keys) {
    $themeColors.Item($slot).RGB = HexToRgbLong $officeColors[$slot]
}
